$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.171.91"
$ws.Range("E2").Value = "  +1.02%  "

# Row 3
$ws.Range("D3").Value = "1.856.57"
$ws.Range("E3").Value = "  +1.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.50%  "

# Row 5
$ws.Range("D5").Value = "'239.12"
$ws.Range("E5").Value = "  +3.65%  "

# Row 6
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.40%  "

# Row 8
$ws.Range("D8").Value = "'42.15"
$ws.Range("E8").Value = "  +7.10%  "

# Row 9
$ws.Range("D9").Value = "'0.330"
$ws.Range("E9").Value = "  +1.42%  "

# Row 10
$ws.Range("E10").Value = "  +1.60%  "

# Row 11
$ws.Range("D11").Value = "'0.0987"
$ws.Range("E11").Value = "  +0.06%  "

# Row 12
$ws.Range("D12").Value = "2.127.04"
$ws.Range("E12").Value = "  +1.52%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.872.44"
$ws.Range("E13").Value = "  +2.27%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.48"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("E15").Value = "  +1.54%  "

# Row 16
$ws.Range("D16").Value = "'4.72"
$ws.Range("E16").Value = "  +1.65%  "

# Row 17
$ws.Range("D17").Value = "35.136.40"
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").Value = "'69.81"
$ws.Range("E18").Value = "  +0.59%  "

# Row 19
$ws.Range("E19").Value = "  +1.17%  "

# Row 20
$ws.Range("D20").Value = "'240.34"
$ws.Range("E20").Value = "  +0.34%  "

# Row 21
$ws.Range("E21").Value = "  +0.35%  "

# Row 22
$ws.Range("D22").Value = "'4.73"
$ws.Range("E22").Value = "  +1.28%  "

# Row 23
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("E24").Value = "  -0.33%  "

# Row 25
$ws.Range("D25").Value = "'168.35"
$ws.Range("E25").Value = "  -2.17%  "

# Row 26
$ws.Range("D26").Value = "'1.92"
$ws.Range("E26").Value = "  +27.75%  "

# Row 27
$ws.Range("E27").Value = "  +3.26%  "

# Row 28
$ws.Range("D28").Value = "'17.63"
$ws.Range("E28").Value = "  +1.87%  "

# Row 29
$ws.Range("E29").Value = "  +0.43%  "

# Row 30
$ws.Range("E30").Value = "  +0.44%  "

# Row 31
$ws.Range("D31").Value = "'0.0558"
$ws.Range("E31").Value = "  +1.43%  "

# Row 32
$ws.Range("D32").Value = "'4.00"
$ws.Range("E32").Value = "  +2.54%  "

# Row 33
$ws.Range("E33").Value = "  +27.36%  "

# Row 34
$ws.Range("E34").Value = "  +2.25%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.03"
$ws.Range("E35").Value = "  +9.90%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.820"
$ws.Range("E36").Value = "  +17.37%  "

# Row 37
$ws.Range("E37").Value = "  +7.06%  "

# Row 38
$ws.Range("E38").Value = "  +4.70%  "

# Row 39
$ws.Range("D39").Value = "'0.0201"
$ws.Range("E39").Value = "  +3.90%  "

# Row 40
$ws.Range("D40").Value = "'89.89"
$ws.Range("E40").Value = "  -1.67%  "

# Row 41
$ws.Range("D41").Value = "1.341.48"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42
$ws.Range("D42").Value = "'0.0591"
$ws.Range("E42").Value = "  +13.44%  "

# Row 43
$ws.Range("D43").Value = "'14.89"
$ws.Range("E43").Value = "  +3.28%  "

# Row 45
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("B46").Value = "Gas"
$ws.Range("C46").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D46").Value = "'12.33"
$ws.Range("E46").Value = "  +44.31%  "

# Row 47
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.74"
$ws.Range("E47").Value = "  -0.42%  "

# Row 48
$ws.Range("D48").Value = "'6.59"
$ws.Range("E48").Value = "  +5.05%  "

# Row 49
$ws.Range("D49").Value = "2.039.45"
$ws.Range("E49").Value = "  +1.37%  "

# Row 50
$ws.Range("D50").Value = "'0.0679"
$ws.Range("E50").Value = "  +1.05%  "

# Row 51
$ws.Range("E51").Value = "  +0.37%  "
